# "exam report and full screen scrolling fix"
# - Replace the 4-question sample quiz with the full 10-question exam bank.
# - Widen option columns C:F so the longer option text is readable.
# - Leave the cursor parked on B9 (mirrors the "full screen scrolling" nav fix).
# - Sheet was saved password-protected ("test"); unprotect before rewriting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("test")

# ---- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Sr"
$ws.Range("B1").Value = "QuestionDesc"
$ws.Range("C1").Value = "Option1"
$ws.Range("D1").Value = "Option2"
$ws.Range("E1").Value = "Option3"
$ws.Range("F1").Value = "Option4"
$ws.Range("G1").Value = "answer"

# ---- Question bank ------------------------------------------------------
$rows = @(
  @(1,  "The 1st person to set foot on Moon was?", "Nil Armstrong", "Rakesh Sharma", "Rayan Prince", "Smith", 1),
  @(2,  "What Mogul emporer of India from 1556-1605, took the throne at age 13 and was the grandson of Babar?", " Delhi Sultanate", "Humayun", "Abbas the Great", "Akbar the Great", 4),
  @(3,  "The well-known theorist of New Social Movement is", "Habermas", "Karl Marx", " Foucault", "Althusser", 1),
  @(4,  "Which was the Napoleon last battle in which he was captured and exiled to St Helena?", "Battle of Waterloo", "Battle of France", "Battle of Paris", "Battle of London", 1),
  @(5,  "In which year Bangladesh was formed as Separate Country?", 1947, 1951, 1956, 1971, 4),
  @(6,  "The longest river in the world is the.", "Nile", "Ganga", "Brhamputra", " Yamuna", 1),
  @(7,  "What country's population had reached an estimated 60 million by the 1570s?", "China", "Japan", " India", "England", 1),
  @(8,  "When was the 1st person to set foot on Moon?", 1970, 1966, 1984, 1977, 2),
  @(9,  "Galileo was an Italian astronomer who:", " Developed the telescope", "Discovered 4 satellites of Jupiter", "Discovered that the movement of the pendulum produces a regular time measurement.", "All are correct", 4),
  @(10, "Which country has always remained free from foreign rule?", " Phillipines", "Nepal", "USA", "Laos", 2)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---- Widen the option columns so the new, longer text fits -------------
$ws.Columns.Item(3).ColumnWidth = 19.833333333333332
$ws.Columns.Item(4).ColumnWidth = 18.5
$ws.Columns.Item(5).ColumnWidth = 14.833333333333334
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668

# ---- Leave the selection on B9, like the saved file -------------------
$ws.Range("B9").Select()
